# "cases index bug fix"
#
# Row 16 of the tracker sheet was mis-entered: the page name said
# "Privacy Notice" and had a "Ждем текст" status highlighted in yellow.
# Fix it:
#   - rename the page to "Privacy Police"
#   - clear the (wrong) status note
#   - drop the yellow highlight, matching the plain formatting used by
#     every other empty status cell in column B (e.g. B3)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the page name in column A.
$ws.Range("A16").Value = "Privacy Police"

# Reset B16's formatting to the plain (no-fill) look used elsewhere in the
# column, by copying an already-plain cell's format over it.
$ws.Range("B3").Copy()
$ws.Range("B16").PasteSpecial(-4122)

# Clear the stale status text now that the highlight is gone.
$ws.Range("B16").Value = ""

# Leave the selection where the user last left it after making the edit.
[void]$ws.Range("B15").Select()
